$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r1 = $ws.Range("B35:AC35")
$r2 = $ws.Range("B36:AC36")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B57:AC57")
$r2 = $ws.Range("B58:AC58")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B164:AC164")
$r2 = $ws.Range("B165:AC165")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B173:AC173")
$r2 = $ws.Range("B174:AC174")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B184:AC184")
$r2 = $ws.Range("B185:AC185")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B186:AC186")
$r2 = $ws.Range("B187:AC187")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B193:AC193")
$r2 = $ws.Range("B194:AC194")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B205:AC205")
$r2 = $ws.Range("B206:AC206")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B207:AC207")
$r2 = $ws.Range("B208:AC208")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B213:AC213")
$r2 = $ws.Range("B214:AC214")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B218:AC218")
$r2 = $ws.Range("B219:AC219")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B283:AC283")
$r2 = $ws.Range("B284:AC284")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B294:AC294")
$r2 = $ws.Range("B295:AC295")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B297:AC297")
$r2 = $ws.Range("B298:AC298")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B311:AC311")
$r2 = $ws.Range("B312:AC312")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B317:AC317")
$r2 = $ws.Range("B318:AC318")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B320:AC320")
$r2 = $ws.Range("B321:AC321")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1

$r1 = $ws.Range("B343:AC343")
$r2 = $ws.Range("B344:AC344")
$v1 = $r1.Value()
$v2 = $r2.Value()
$r1.Value = $v2
$r2.Value = $v1
